$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table was re-laid out into a new column arrangement:
#   - country/territory names moved from column B into column A
#   - the "Country/Territory" / "Subtotal" / "TOTAL" labels that used to sit in
#     column C were dropped entirely
#   - the week-number labels moved from column D into column B
#   - the data values moved from column E into column D
#   - columns F, G and H (incidence rate / deaths / population) stayed put
# Clear the old, now-obsolete layout first.
$ws.Range("A1:I38").ClearContents()

# Write the values into their new locations
$ws.Range("H1").Value = 'Population d'
$ws.Range("B2").Value = 'Weeka'
$ws.Range("D2").Value = 'Confirmed cases b'
$ws.Range("F2").Value = 'Incidence ratec'
$ws.Range("G2").Value = 'Deaths'
$ws.Range("H2").Value = 'X 1000'
$ws.Range("A3").Value = 'Latin Caribbean'
$ws.Range("A4").Value = 'Cuba'
$ws.Range("D4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = '''11,266'
$ws.Range("A5").Value = 'Dominican Republic'
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = '''10,404'
$ws.Range("A6").Value = 'French Guiana (*)'
$ws.Range("B6").Value = 'Week 8'
$ws.Range("D6").Value = 7
$ws.Range("F6").Value = 2.8
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 249
$ws.Range("A7").Value = 'Guadaloupe (**)'
$ws.Range("B7").Value = 'Week 8'
$ws.Range("D7").Value = 335
$ws.Range("F7").Value = 71.9
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 466
$ws.Range("A8").Value = 'Haiti'
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = '''10,317'
$ws.Range("A9").Value = 'Martinique'
$ws.Range("B9").Value = 'Week 8'
$ws.Range("D9").Value = 943
$ws.Range("F9").Value = 233.4
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 404
$ws.Range("A10").Value = 'Puerto Rico'
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = '''3,688'
$ws.Range("A11").Value = 'Saint Barthelemy'
$ws.Range("B11").Value = 'Week 8'
$ws.Range("D11").Value = 114
$ws.Range("F11").Value = '''1,280.9'
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 9
$ws.Range("A12").Value = 'Saint Martin (French part) (***)'
$ws.Range("B12").Value = 'Week 8'
$ws.Range("D12").Value = 711
$ws.Range("F12").Value = '''1,992.2'
$ws.Range("G12").Value = '1#'
$ws.Range("H12").Value = 36
$ws.Range("D13").Value = '''2,110'
$ws.Range("F13").Value = 5.7
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = '''36,839'
$ws.Range("A14").Value = 'Non-Latin Caribbean'
$ws.Range("A15").Value = 'Anguilla$'
$ws.Range("B15").Value = 'Week 8'
$ws.Range("D15").Value = 11
$ws.Range("F15").Value = 68.8
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 16
$ws.Range("A16").Value = 'Antigua & Barbuda'
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 90
$ws.Range("A17").Value = 'Aruba$'
$ws.Range("B17").Value = 'Week 6'
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 0.9
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 109
$ws.Range("A18").Value = 'Bahamas'
$ws.Range("D18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 377
$ws.Range("A19").Value = 'Barbados'
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 289
$ws.Range("A20").Value = 'Cayman Islands'
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 54
$ws.Range("A21").Value = 'Curacao'
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 147
$ws.Range("A22").Value = 'Dominica (****)'
$ws.Range("B22").Value = 'Week 7'
$ws.Range("D22").Value = 45
$ws.Range("F22").Value = 61.6
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 73
$ws.Range("A23").Value = 'Grenada'
$ws.Range("D23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 110
$ws.Range("A24").Value = 'Guyana'
$ws.Range("D24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 800
$ws.Range("A25").Value = 'Jamaica'
$ws.Range("D25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = '''2,784'
$ws.Range("A26").Value = 'Montserrat'
$ws.Range("D26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 5
$ws.Range("A27").Value = 'Saint Kitts & Nevis'
$ws.Range("B27").Value = 'Week 8'
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 51
$ws.Range("A28").Value = 'Saint Lucia'
$ws.Range("D28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 163
$ws.Range("A29").Value = 'Saint Vincent & the Grenadines'
$ws.Range("D29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 103
$ws.Range("A30").Value = 'Sint Maarten (Dutch part)'
$ws.Range("B30").Value = 'Week 6'
$ws.Range("D30").Value = 65
$ws.Range("F30").Value = 162.5
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 40
$ws.Range("A31").Value = 'Suriname'
$ws.Range("D31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 539
$ws.Range("A32").Value = 'Trinidad & Tobago'
$ws.Range("D32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = '''1,341'
$ws.Range("A33").Value = 'Turks & Caicos Islands'
$ws.Range("D33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 48
$ws.Range("A34").Value = 'Virgin Islands (UK)'
$ws.Range("B34").Value = 'Week 5'
$ws.Range("D34").Value = 5
$ws.Range("F34").Value = 15.6
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 32
$ws.Range("A35").Value = 'Virgin Islands (US)'
$ws.Range("D35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 105
$ws.Range("D36").Value = 128
$ws.Range("F36").Value = 1.8
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = '''7,276'
$ws.Range("D37").Value = '''2,238'
$ws.Range("F37").Value = 5.1
$ws.Range("G37").Value = '-'
$ws.Range("H37").Value = '''44,115'
$ws.Range("A38").Value = 'NOTES'

# Remove the text-as-number formatting marker Excel applied when the numeric-
# looking values above (e.g. "11,266") were entered as text, so no stray cell
# styling is left behind.
$ws.Cells.ClearFormats()
